$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.304492712020874
$ws.Range("B1").Value = 2.663942813873291
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.645938396453857
$ws.Range("E1").Value = 1.100188851356506
